$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsDesc = $wb.Worksheets.Item("Description")

# --- Description sheet: update the "Source" note into a proper
# Source / References block, per PM&C-supplied text. ---

# Row 9: turn the old single "Sourced from: ..." note into a "Source" label
# (col A) + first reference (col B), matching the style already used for
# the A3:B8 label/body pairs.
$wsDesc.Range("A9").Value = "Source"
$wsDesc.Range("B9").Value = "ABS (unpublished) Australian Health Survey 2014–15"

# Rows 10-11: additional source lines under the same "Source" label.
$wsDesc.Range("B10").Value = "ABS (unpublished) Australian Health Survey, 2011-13 (2011-12 Core component)"
$wsDesc.Range("B11").Value = "ABS (unpublished), National Health Survey 2007-08."

# Row 12: new "References" label + citation.
$wsDesc.Range("A12").Value = "References"
$wsDesc.Range("B12").Value = "Haire-Joshu, D. and Nanney, M., 2002, Prevention of Overweight and Obesity in Children: Influences on the Food Environment. The Diabetes Educator, 28(3), pp.415-423."

# Match formatting of the existing label / body columns.
$wsDesc.Range("A9:A11").Style = $wsDesc.Range("A8").Style
$wsDesc.Range("B9:B11").Style = $wsDesc.Range("B5").Style
$wsDesc.Range("B9:B11").WrapText = $true

# "References" row sits outside the label/body pattern used above it -
# plain default cell for the label, larger 12pt body text for the citation.
$wsDesc.Range("B12").Font.Size = 12
$wsDesc.Range("B12").WrapText = $true

$wsDesc.Rows.Item(9).RowHeight = 13.8
$wsDesc.Rows.Item(10).RowHeight = 13.8
$wsDesc.Rows.Item(11).RowHeight = 13.8
$wsDesc.Rows.Item(12).RowHeight = 26.95

# --- Selection / active-sheet bookkeeping left over from the edit session ---
[void]$wsData.Range("B9:B14").Select()
[void]$wsDesc.Range("B9:B14").Select()
[void]$wsDesc.Activate()
